$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell F1, styled like the other header cells (B1:E1)
$ws.Range("F1").Value = "time_taken"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)   # xlPasteFormats: reuse E1's header style

# Timestamp values for rows 2-14 (as text, same as the other inlineStr data cells)
$timestamps = @(
    "2021-10-05 10:50:43.169481",
    "2021-10-05 10:50:43.169494",
    "2021-10-05 10:50:43.169498",
    "2021-10-05 10:50:43.169501",
    "2021-10-05 10:50:43.169505",
    "2021-10-05 10:50:43.169508",
    "2021-10-05 10:50:43.169512",
    "2021-10-05 10:50:43.169515",
    "2021-10-05 10:50:43.169518",
    "2021-10-05 10:50:43.169521",
    "2021-10-05 10:50:43.169525",
    "2021-10-05 10:50:43.169528",
    "2021-10-05 10:50:43.169531"
)

for ($i = 0; $i -lt $timestamps.Length; $i++) {
    $row = $i + 2
    $cell = $ws.Cells.Item($row, 6)
    $cell.Value = $timestamps[$i]
}
